$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.649.84"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "1.867.41"

$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'234.84"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.4727"
$ws.Range("E7").Value = "  -0.80%  "

$ws.Range("D8").Value = "'0.2769"
$ws.Range("E8").Value = "  +0.46%  "

$ws.Range("D9").Value = "'0.06372"
$ws.Range("E9").Value = "  -1.33%  "

$ws.Range("D10").Value = "'17.77"
$ws.Range("E10").Value = "  +9.73%  "

$ws.Range("D11").Value = "1.867.88"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").Value = "'0.07471"
$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("D13").Value = "'4.981"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").Value = "'85.26"
$ws.Range("E14").Value = "  -0.83%  "

$ws.Range("D15").Value = "'0.6322"
$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("D16").Value = "30.591.42"
$ws.Range("E16").Value = "  +0.93%  "

$ws.Range("D17").Value = "'241.41"
$ws.Range("E17").Value = "  +4.52%  "

$ws.Range("D18").Value = "'0.9989"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").Value = "'12.77"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("D20").Value = "'0.000007373"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").Value = "'4.998"
$ws.Range("E22").Value = "  -2.10%  "

$ws.Range("D23").Value = "'5.985"
$ws.Range("E23").Value = "  -0.78%  "

$ws.Range("D24").Value = "'9.374"
$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("D25").Value = "'166.32"
$ws.Range("E25").Value = "  -0.66%  "

$ws.Range("D26").Value = "'18.19"
$ws.Range("E26").Value = "  +1.48%  "

$ws.Range("D27").Value = "'1.884"
$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("D28").Value = "'0.1031"
$ws.Range("E28").Value = "  +2.81%  "

$ws.Range("D29").Value = "'1.381"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").Value = "'4.108"
$ws.Range("E30").Value = "  -2.90%  "

$ws.Range("D31").Value = "'3.867"
$ws.Range("E31").Value = "  -1.44%  "

$ws.Range("D32").Value = "'0.04933"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").Value = "'1.152"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").Value = "'0.7086"
$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("D35").Value = "'2.704"
$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").Value = "'0.01905"
$ws.Range("E36").Value = "  -1.66%  "

$ws.Range("D37").Value = "'2.692"
$ws.Range("E37").Value = "  +2.17%  "

$ws.Range("D38").Value = "'0.8835"
$ws.Range("E38").Value = "  -2.67%  "

$ws.Range("D39").Value = "'2.001"
$ws.Range("E39").Value = "  +0.52%  "

$ws.Range("D40").Value = "'105.82"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("D41").Value = "'0.9998"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "'0.4108"
$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("D43").Value = "'5.553"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").Value = "'7.243"
$ws.Range("E44").Value = "  +2.22%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'62.25"
$ws.Range("E45").Value = "  +1.44%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1232"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("D47").Value = "'8.661"
$ws.Range("E47").Value = "  -1.16%  "

$ws.Range("D48").Value = "'33.74"
$ws.Range("E48").Value = "  +2.01%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05576"
$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.385"
$ws.Range("E50").Value = "  -1.30%  "

$ws.Range("D51").Value = "'0.3704"
$ws.Range("E51").Value = "  -0.39%  "
